$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $ok = $f.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "MISSING: $findText"
    }
    return $ok
}

# --- 1) Title line: "TRISTAR AERO TECHNOLOGY, LLC" + trailing two spaces -> "TRISTAR AERO TECHNOLOGY, " + bold/italic "INC"
$rngTitle = $d.Content
$fTitle = $rngTitle.Find
$fTitle.Execute("TRISTAR AERO TECHNOLOGY, LLC  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($fTitle.Found) {
    $titleStart = $rngTitle.Start
    $titleEnd = $rngTitle.End
    # Replace whole match (company name + LLC + two trailing spaces) with new text in one go,
    # keeping the big bold/italic formatting that the original run already has.
    $rngTitle.Text = "TRISTAR AERO TECHNOLOGY, INC"
} else {
    Write-Output "MISSING: title line"
}

# --- 2) Pricing table cells ---
ReplaceText "＄600.00 USD" "＄650.00 USD"
ReplaceText "＄550.00 USD(二個使用者共＄1100.00 USD)  " "＄620.00 USD(二個使用者共＄1240.00 USD)  "
ReplaceText "＄500.00 USD(三個使用者共＄1500.00 USD 每加一個使用者加＄500.00 USD)  " "＄590.00 USD(三個使用者共＄1770.00 USD 每加一個使用者加＄590.00 USD)  "
ReplaceText "＄1100.00 USD    " "＄1240.00 USD    "
ReplaceText "＄1500.00 USD  " "＄1770.00 USD  "
ReplaceText "＄2000.00 USD " "＄2360.00 USD "
ReplaceText "＄2700.00 USD " "＄3360.00 USD "
ReplaceText "＄2700.00 USD (每加一個使用者加＄900.00 USD )  " "＄3540.00 USD (每加一個使用者加＄1180.00 USD )  "
ReplaceText "＄3600.00 USD (每加一個使用者加＄1200.00 USD )  " "＄4770.00 USD (每加一個使用者加＄1590.00 USD )  "

# --- 3) "(TRISTAR AERO TECHNOLOGY, LLC  " -> "(TRISTAR AERO TECHNOLOGY, INC" + two spaces (LLC -> INC, merged run)
ReplaceText "(TRISTAR AERO TECHNOLOGY, LLC  " "(TRISTAR AERO TECHNOLOGY, INC  "

Write-Output "done"
